# Update "想去人数" (interested-people count) values across the sheets
# to reflect the newly generated data snapshot (gh-pages output at 456a3b4).

$wb = $excel.ActiveWorkbook

# --- Sheet: 展览 (Exhibitions) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value  = 240
$ws1.Range("F5").Value  = 1900
$ws1.Range("F6").Value  = 218
$ws1.Range("F7").Value  = 642
$ws1.Range("F8").Value  = 13
$ws1.Range("F9").Value  = 130
$ws1.Range("F11").Value = 626
$ws1.Range("F12").Value = 19
$ws1.Range("F13").Value = 72
$ws1.Range("F14").Value = 591

# --- Sheet: 演出 (Performances) ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F6").Value  = 7
$ws2.Range("F11").Value = 24
$ws2.Range("F12").Value = 206

# --- Sheet: 本地生活 (Local life) ---
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F2").Value = 6309
$ws3.Range("F4").Value = 1959
$ws3.Range("F5").Value = 172

# --- Sheet: 全部类型 (All types - combined view) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value  = 6309
$ws4.Range("F4").Value  = 1959
$ws4.Range("F6").Value  = 172
$ws4.Range("F12").Value = 240
$ws4.Range("F13").Value = 7
$ws4.Range("F16").Value = 1900
$ws4.Range("F18").Value = 218
$ws4.Range("F20").Value = 24
$ws4.Range("F21").Value = 642
$ws4.Range("F22").Value = 13
$ws4.Range("F23").Value = 130
$ws4.Range("F24").Value = 206
$ws4.Range("F26").Value = 626
$ws4.Range("F27").Value = 19
$ws4.Range("F28").Value = 72
$ws4.Range("F30").Value = 592
